# The edit re-orders four observation records (rows 5-8): the data that
# used to live in row 5 ends up in row 7 and vice versa, and likewise the
# data that used to live in row 6 ends up in row 8 and vice versa. All
# other rows are left untouched.
#
# We implement this as two full-row swaps using Copy / PasteSpecial
# (values), staging through a scratch row far below the used range so
# that genuinely blank source cells correctly blank out the destination
# (PasteSpecial leaves a cell alone if the corresponding source cell is
# blank, so the destination must be cleared first). Copy/PasteSpecial is
# used instead of direct Value assignment so that text that looks like a
# date ("2020-06-11") or time ("00:00") is carried over as plain text
# instead of being re-interpreted/reformatted by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = "AY"
$scratchRow = 1000
$scratchRange = "A" + $scratchRow + ":" + $lastCol + $scratchRow
$xlPasteValues = -4163

function Swap-EntireRows {
    param(
        [int]$RowA,
        [int]$RowB
    )

    $rangeA = "A" + $RowA + ":" + $lastCol + $RowA
    $rangeB = "A" + $RowB + ":" + $lastCol + $RowB

    # scratch = RowA
    $ws.Range($scratchRange).Clear()
    $ws.Range($rangeA).Copy()
    $ws.Range($scratchRange).PasteSpecial($xlPasteValues)

    # RowA = RowB
    $ws.Range($rangeA).Clear()
    $ws.Range($rangeB).Copy()
    $ws.Range($rangeA).PasteSpecial($xlPasteValues)

    # RowB = scratch (original RowA)
    $ws.Range($rangeB).Clear()
    $ws.Range($scratchRange).Copy()
    $ws.Range($rangeB).PasteSpecial($xlPasteValues)

    $ws.Range($scratchRange).Clear()
}

Swap-EntireRows 5 7
Swap-EntireRows 6 8

$excel.CutCopyMode = $false
